$p = $ppt.ActivePresentation
$cs = $p.ColorSchemes
Write-Output $cs.Count
for ($i=1; $i -le $cs.Count; $i++) {
  try {
    $item = $cs.Item($i)
    Write-Output "$i -> $item"
  } catch {
    Write-Output "$i -> ERR $_"
  }
}
